$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (volume/number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -66.666666666666
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 42.857142857142
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 233.333333333333
$ws.Range("N15").Value = -47.368421052631
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 84
$ws.Range("K16").Value = 11.904761904761
$ws.Range("L16").Value = 22.077922077922
$ws.Range("M16").Value = -36.486486486486
$ws.Range("N16").Value = -82.90909090909
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 29.166666666666
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = 2.189781021897
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 64.705882352941
$ws.Range("N17").Value = -38.59649122807
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 21.052631578947
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = 68.421052631578
$ws.Range("L18").Value = 12.676056338028
$ws.Range("M18").Value = -6.432748538011
$ws.Range("N18").Value = -72.076788830715
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = -71.428571428571
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 319
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = -4.491017964071
$ws.Range("L19").Value = 18.148148148148
$ws.Range("M19").Value = 69.680851063829
$ws.Range("N19").Value = 37.5
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 61
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = -14.084507042253
$ws.Range("L20").Value = -3.174603174603
$ws.Range("M20").Value = -12.857142857142
$ws.Range("N20").Value = -83.947368421052
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -27.659574468085
$ws.Range("F21").Value = 148
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = 8.823529411764
$ws.Range("I21").Value = 787
$ws.Range("J21").Value = 732
$ws.Range("K21").Value = 7.513661202185
$ws.Range("L21").Value = 16.076696165191
$ws.Range("M21").Value = 18.168168168168
$ws.Range("N21").Value = -60.472124560522
$ws.Range("F22").Value = 1
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = -45.454545454545
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -14.285714285714
$ws.Range("I23").Value = 84
$ws.Range("J23").Value = 89
$ws.Range("K23").Value = -5.617977528089
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = 40
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = 10.679611650485
$ws.Range("I24").Value = 507
$ws.Range("J24").Value = 480
$ws.Range("K24").Value = 5.625
$ws.Range("L24").Value = -9.464285714285
$ws.Range("M24").Value = -4.158790170132
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 171
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 185
$ws.Range("L25").Value = 36.8
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 23.255813953488
$ws.Range("I26").Value = 267
$ws.Range("J26").Value = 208
$ws.Range("K26").Value = 28.365384615384
$ws.Range("L26").Value = 16.086956521739
$ws.Range("M26").Value = 19.730941704035
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 13
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 44.444444444444
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -13.333333333333
$ws.Range("L28").Value = 23.809523809523
$ws.Range("N29").Value = -81.25
$ws.Range("N30").Value = -83.870967741935
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 15
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = 36.363636363636

# --- Convert numeric cells to text placeholders ("0" / "***.*"), style -> 14 ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("G14,H14,D29,E29,D30,E30").PasteSpecial(-4122)

# --- Convert text placeholder cells to numbers, style -> 15 (count columns) ---
$ws.Range("C15").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("C27").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C15,D22,G22,C27,D28,D31,G31").PasteSpecial(-4122)

# --- Convert text placeholder cells to numbers, style -> 16 (pct-chg columns) ---
$ws.Range("E22").Value = -100
$ws.Range("H22").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("E31").Value = -100
$ws.Range("H31").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E22,H22,E28,E31,H31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
